# fix!: rename sheet CONDUCTOR_COUPLING.
$wb = $excel.ActiveWorkbook

# Rename the "CONDUCTOR_COUPLING" sheet to "CONDUCTOR_coupling"
$ws = $wb.Worksheets.Item("CONDUCTOR_COUPLING")
$ws.Name = "CONDUCTOR_coupling"

# Make it the active/selected sheet (was CONDUCTOR_operation before)
$ws.Activate()
